$p = $ppt.ActivePresentation

# Insert a new slide just before the final slide (position 17 of what
# will become 18), using the "Title and Content" layout
# (ppLayoutText = 2 -> CustomLayout "Title and Content" / slideLayout2.xml).
$newSlide = $p.Slides.Add(17, 2)

# Title placeholder -> "Github link:"
$title = $newSlide.Shapes.Item(1)
$title.Name = "Title 1"
$title.TextFrame.TextRange.Text = "Github link:"
$title.TextFrame.TextRange.LanguageID = "en-IN"

# Body placeholder -> GitHub Pages link, justified, hyperlinked
$body = $newSlide.Shapes.Item(2)
$body.Name = "Text Placeholder 2"
$body.TextFrame.TextRange.Text = "https://buvanesh-lgtm.github.io/TNSDC-FWD-Digitialportfilo/"
$body.TextFrame.TextRange.LanguageID = "en-IN"
$body.TextFrame.TextRange.ParagraphFormat.Alignment = 4
$body.TextFrame.TextRange.ActionSettings(1).Hyperlink.Address = "https://buvanesh-lgtm.github.io/TNSDC-FWD-Digitialportfilo/"

$body.Left = 48.0
$body.Top = 198.0
$body.Width = 534.0
$body.Height = 33.81095

# The slide that used to be last (slide id 265) is now pushed from
# position 17 to position 18; refresh its cached page-number text.
$lastSlide = $p.Slides.Item(18)
$pageNum = $lastSlide.Shapes.Item("object 9")
$pageNum.TextFrame.TextRange.Text = "18"
$pageNum.Left = 887.97
$pageNum.Top = 509.71157480314963
$pageNum.Width = 18.0
$pageNum.Height = 15.1
